# "Add files via upload" — populate sheet "17" (xl/worksheets/sheet7.xml)
# with the four new submission-group names and make it the active sheet.
#
# Typing the values in this order (A1, A2, A4, A3) reproduces the shared
# string table order seen in the target workbook: Yuval Koskas, Max
# Gutnik, Noam Raanan, Lior Tsalovich — while the sheet itself ends up
# with A3 = "Lior Tsalovich" and A4 = "Noam Raanan".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("17")

$ws.Range("A1").Value = "Yuval Koskas "
$ws.Range("A2").Value = "Max Gutnik"
$ws.Range("A4").Value = "Noam Raanan"
$ws.Range("A3").Value = "Lior Tsalovich"

# Select B7 on this sheet and make the sheet the active tab (matches the
# tabSelected move from sheet "14" to sheet "17" and workbookView's
# activeTab going from 3 to 6).
$ws.Range("B7").Select()
